# Apply updated crypto market data (price + 1h volume change) per GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text like "54.558.73" / "1.00" — force Text format first so
# Excel does not reinterpret these digit-and-dot strings as numbers/dates.
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D8", "D9", "D10", "D12", "D13", "D14", "D15", "D16", "D18", "D19", "D21", "D22", "D23", "D24", "D25", "D28", "D30", "D31", "D32", "D34", "D35", "D36", "D39", "D42", "D43", "D44", "D45", "D49")
foreach ($c in $priceCells) { $ws.Range($c).NumberFormat = "@" }

$ws.Range("D2").Value = "54.558.73"
$ws.Range("E2").Value = "  +0.32%  "

$ws.Range("D3").Value = "2.282.09"
$ws.Range("E3").Value = "  -0.08%  "

$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").Value = "505.10"
$ws.Range("E5").Value = "  +1.07%  "

$ws.Range("D6").Value = "128.60"
$ws.Range("E6").Value = "  -0.47%  "

$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "0.529"
$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").Value = "2.299.51"
$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("D10").Value = "0.0967"
$ws.Range("E10").Value = "  +1.23%  "

$ws.Range("E11").Value = "  +1.69%  "

$ws.Range("D12").Value = "0.341"
$ws.Range("E12").Value = "  +1.99%  "

$ws.Range("D13").Value = "4.90"
$ws.Range("E13").Value = "  +3.74%  "

$ws.Range("D14").Value = "23.43"
$ws.Range("E14").Value = "  +1.63%  "

$ws.Range("D15").Value = "2.691.11"
$ws.Range("E15").Value = "  -0.01%  "

$ws.Range("D16").Value = "54.617.26"
$ws.Range("E16").Value = "  +0.56%  "

$ws.Range("E17").Value = "  +1.04%  "

$ws.Range("D18").Value = "2.282.60"
$ws.Range("E18").Value = "  +0.30%  "

$ws.Range("D19").Value = "10.40"
$ws.Range("E19").Value = "  +1.33%  "

$ws.Range("E20").Value = "  +0.47%  "

$ws.Range("D21").Value = "306.43"
$ws.Range("E21").Value = "  +0.32%  "

$ws.Range("D22").Value = "6.55"
$ws.Range("E22").Value = "  +2.23%  "

$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").Value = "60.38"
$ws.Range("E24").Value = "  -2.67%  "

$ws.Range("D25").Value = "0.996"
$ws.Range("E25").Value = "  -0.20%  "

$ws.Range("E26").Value = "  -1.05%  "

$ws.Range("E27").Value = "  +1.74%  "

$ws.Range("D28").Value = "171.36"
$ws.Range("E28").Value = "  -1.74%  "

$ws.Range("E29").Value = "  +1.64%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "1.63"
$ws.Range("E30").Value = "  +0.74%  "

$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0702"
$ws.Range("E31").Value = "  +1.68%  "

$ws.Range("D32").Value = "1.13"
$ws.Range("E32").Value = "  +4.22%  "

$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("D34").Value = "17.97"
$ws.Range("E34").Value = "  +0.80%  "

$ws.Range("D35").Value = "0.994"
$ws.Range("E35").Value = "  -0.26%  "

$ws.Range("D36").Value = "0.908"
$ws.Range("E36").Value = "  -2.82%  "

$ws.Range("E37").Value = "  +0.17%  "

$ws.Range("E38").Value = "  +1.02%  "

$ws.Range("D39").Value = "36.55"
$ws.Range("E39").Value = "  +1.37%  "

$ws.Range("E40").Value = "  +0.24%  "

$ws.Range("E41").Value = "  +0.37%  "

$ws.Range("D42").Value = "5.05"
$ws.Range("E42").Value = "  +3.75%  "

$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "3.40"
$ws.Range("E43").Value = "  +0.15%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "130.97"
$ws.Range("E44").Value = "  +4.71%  "

$ws.Range("D45").Value = "251.21"
$ws.Range("E45").Value = "  +4.06%  "

$ws.Range("E46").Value = "  +0.49%  "

$ws.Range("E47").Value = "  +1.47%  "

$ws.Range("E48").Value = "  +0.39%  "

$ws.Range("D49").Value = "0.375"
$ws.Range("E49").Value = "  +0.53%  "

$ws.Range("E50").Value = "  +0.35%  "

$ws.Range("E51").Value = "  +0.41%  "
